$d = $word.ActiveDocument

# Locate the date "14.05.2024" that needs correcting to "13.05.2024" inside
# the "Drugie spotkanie organizacyjne" paragraph (search with enough context
# so we only match this specific occurrence of the date).
$findRng = $d.Content
$findRng.Find.Execute("do 14.05.2024. Jako", $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)

if ($findRng.Find.Found) {
    # Within the matched range, "do 1" is 4 characters, so the "4" digit that
    # needs to become "3" sits right after that.
    $digitStart = $findRng.Start + 4
    $digitEnd = $digitStart + 1
    $digitRng = $d.Range($digitStart, $digitEnd)

    # Replace just the single "4" with "3".
    $digitRng.Text = "3"

    # Re-apply the formatted text to the edited range so it settles into its
    # own run boundary (matching how Word keeps the freshly typed character
    # in a distinct run from its unmodified neighbours).
    $digitRng.FormattedText = $digitRng.FormattedText
}
